# Reorder / rename the parameter columns and fix up the
# "verbose_list.no_defaults" row so the pandas reference parameters sheet
# matches the parameter library's current column naming (expression,
# maximum, minimum, non_negative, standard_error, value, vary).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$header = New-Object 'object[,]' 1,8
$header[0,0] = "label"
$header[0,1] = "expression"
$header[0,2] = "maximum"
$header[0,3] = "minimum"
$header[0,4] = "non_negative"
$header[0,5] = "standard_error"
$header[0,6] = "value"
$header[0,7] = "vary"
$ws.Range("A1:H1").Value = $header

# --- Data rows --------------------------------------------------------
# Columns (new order): A label | B expression | C maximum | D minimum |
#                       E non_negative | F standard_error | G value | H vary

# row 2 - pure_list.1
$ws.Cells.Item(2,2).Value = "None"
$ws.Cells.Item(2,3).Value = ""
$ws.Cells.Item(2,4).Value = ""
$ws.Cells.Item(2,5).Value = $false
$ws.Cells.Item(2,6).Value = "None"
$ws.Cells.Item(2,7).Value = 1
$ws.Cells.Item(2,8).Value = $true

# row 3 - pure_list.2
$ws.Cells.Item(3,2).Value = "None"
$ws.Cells.Item(3,3).Value = ""
$ws.Cells.Item(3,4).Value = ""
$ws.Cells.Item(3,5).Value = $false
$ws.Cells.Item(3,6).Value = "None"
$ws.Cells.Item(3,7).Value = 2
$ws.Cells.Item(3,8).Value = $true

# row 4 - list_with_options.1
$ws.Cells.Item(4,2).Value = "None"
$ws.Cells.Item(4,3).Value = ""
$ws.Cells.Item(4,4).Value = ""
$ws.Cells.Item(4,5).Value = $false
$ws.Cells.Item(4,6).Value = "None"
$ws.Cells.Item(4,7).Value = 3
$ws.Cells.Item(4,8).Value = $false

# row 5 - list_with_options.2
$ws.Cells.Item(5,2).Value = "None"
$ws.Cells.Item(5,3).Value = ""
$ws.Cells.Item(5,4).Value = ""
$ws.Cells.Item(5,5).Value = $false
$ws.Cells.Item(5,6).Value = "None"
$ws.Cells.Item(5,7).Value = 4
$ws.Cells.Item(5,8).Value = $false

# row 6 - verbose_list.all_defaults
$ws.Cells.Item(6,2).Value = "None"
$ws.Cells.Item(6,3).Value = ""
$ws.Cells.Item(6,4).Value = ""
$ws.Cells.Item(6,5).Value = $false
$ws.Cells.Item(6,6).Value = "None"
$ws.Cells.Item(6,7).Value = 5
$ws.Cells.Item(6,8).Value = $true

# row 7 - verbose_list.no_defaults (minimum/maximum now populated)
$ws.Cells.Item(7,2).Value = "None"
$ws.Cells.Item(7,3).Value = 1
$ws.Cells.Item(7,4).Value = -1
$ws.Cells.Item(7,5).Value = $true
$ws.Cells.Item(7,6).Value = "None"
$ws.Cells.Item(7,7).Value = 6
$ws.Cells.Item(7,8).Value = $false

# row 8 - verbose_list.expression_only
$ws.Cells.Item(8,2).Value = "`$verbose_list.all_defaults + `$verbose_list.no_defaults"
$ws.Cells.Item(8,3).Value = ""
$ws.Cells.Item(8,4).Value = ""
$ws.Cells.Item(8,5).Value = $false
$ws.Cells.Item(8,6).Value = "None"
$ws.Cells.Item(8,7).Value = 11
$ws.Cells.Item(8,8).Value = $false
